$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 header: reuse the same formatting as the other header cells (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:52:32.420023",
    "2021-10-05 10:52:32.420033",
    "2021-10-05 10:52:32.420037",
    "2021-10-05 10:52:32.420039",
    "2021-10-05 10:52:32.420042",
    "2021-10-05 10:52:32.420045",
    "2021-10-05 10:52:32.420048",
    "2021-10-05 10:52:32.420050",
    "2021-10-05 10:52:32.420053",
    "2021-10-05 10:52:32.420056"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
